$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 241.85715
$ws.Range("I12").Value = 249.66667
$ws.Range("J12").Value = 195
$ws.Range("K12").Value = 249.66667
$ws.Range("L12").Value = 195
$ws.Range("M12").Value = -79.66667000000001
$ws.Range("N12").Value = -535

$ws.Range("H40").Value = 1332.8182
$ws.Range("I40").Value = 1325.8572
$ws.Range("J40").Value = 1345
$ws.Range("K40").Value = 1325.8572
$ws.Range("L40").Value = 1345
$ws.Range("M40").Value = -1150.8572
$ws.Range("N40").Value = -1695

$ws.Range("H98").Value = 698
$ws.Range("I98").Value = 682.3077
$ws.Range("J98").Value = 800
$ws.Range("K98").Value = 682.3077
$ws.Range("L98").Value = 800
$ws.Range("M98").Value = 815.6923
$ws.Range("N98").Value = -3796

$ws.Range("H122").Value = 698
$ws.Range("I122").Value = 682.3077
$ws.Range("J122").Value = 800
$ws.Range("K122").Value = 2046.9231
$ws.Range("L122").Value = 2400
$ws.Range("M122").Value = 403.0769
$ws.Range("N122").Value = -7300

$ws.Range("H125").Value = 1859.7778
$ws.Range("I125").Value = 810.6667
$ws.Range("J125").Value = 2384.3333
$ws.Range("K125").Value = 7296.0003
$ws.Range("L125").Value = 21458.9997
$ws.Range("M125").Value = -4836.0003
$ws.Range("N125").Value = -26378.9997

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H9").Value = 19920
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 19920
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 19920
$ws.Range("N9").Value = -20260

$ws.Range("H20").Value = 19920
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 19920
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 19920
$ws.Range("N20").Value = -20460

$ws.Range("H23").Value = 13832.5
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 13832.5
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 13832.5
$ws.Range("N23").Value = -14350.5

$ws.Range("H36").Value = 13333.333
$ws.Range("I36").Value = 10000
$ws.Range("J36").Value = 15000
$ws.Range("K36").Value = 10000
$ws.Range("L36").Value = 15000
$ws.Range("M36").Value = -9654
$ws.Range("N36").Value = -15692

$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()

$ws.Range("H88").Value = 2816.25
$ws.Range("I88").Value = 2624.6
$ws.Range("J88").Value = 3135.6667
$ws.Range("K88").Value = 2624.6
$ws.Range("L88").Value = 3135.6667
$ws.Range("M88").Value = -2218.6
$ws.Range("N88").Value = -3947.6667

$ws.Range("H91").Value = 2816.25
$ws.Range("I91").Value = 2624.6
$ws.Range("J91").Value = 3135.6667
$ws.Range("K91").Value = 2624.6
$ws.Range("L91").Value = 3135.6667
$ws.Range("M91").Value = -1220.6
$ws.Range("N91").Value = -5943.6667

$ws.Range("H122").Value = 1462.2
$ws.Range("I122").Value = 1312
$ws.Range("J122").Value = 1499.75
$ws.Range("K122").Value = 3936
$ws.Range("L122").Value = 4499.25
$ws.Range("M122").Value = -1486
$ws.Range("N122").Value = -9399.25

$ws.Range("H132").Value = 2554.5945
$ws.Range("I132").Value = 2285.1724
$ws.Range("J132").Value = 3531.25
$ws.Range("K132").Value = 6855.5172
$ws.Range("L132").Value = 10593.75
$ws.Range("M132").Value = -4325.5172
$ws.Range("N132").Value = -15653.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H25").Value = 1401.7142
$ws.Range("I25").Value = 1067.6364
$ws.Range("J25").Value = 2626.6667
$ws.Range("K25").Value = 1067.6364
$ws.Range("L25").Value = 2626.6667
$ws.Range("M25").Value = -832.6364000000001
$ws.Range("N25").Value = -3096.6667

$ws.Range("H54").Value = 496.66666
$ws.Range("I54").Value = 496.66666
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 496.66666
$ws.Range("L54").Value = 0
$ws.Range("M54").Value = -12.66665999999998
$ws.Range("N54").ClearContents()

$ws.Range("H56").Value = 12000
$ws.Range("I56").Value = 0
$ws.Range("J56").Value = 12000
$ws.Range("K56").Value = 0
$ws.Range("L56").Value = 12000
$ws.Range("N56").Value = -13478

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 3035.55
$ws.Range("I99").Value = 2421.2
$ws.Range("J99").Value = 3649.9
$ws.Range("K99").Value = 2421.2
$ws.Range("L99").Value = 3649.9
$ws.Range("M99").Value = -923.1999999999998
$ws.Range("N99").Value = -6645.9

$ws.Range("H107").Value = 649.43475
$ws.Range("I107").Value = 584.8461
$ws.Range("J107").Value = 733.4
$ws.Range("K107").Value = 584.8461
$ws.Range("L107").Value = 733.4
$ws.Range("M107").Value = 1335.1539
$ws.Range("N107").Value = -4573.4

$ws.Range("H122").Value = 996
$ws.Range("I122").Value = 996
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2988
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -538
$ws.Range("N122").ClearContents()

$ws.Range("H126").Value = 3035.55
$ws.Range("I126").Value = 2421.2
$ws.Range("J126").Value = 3649.9
$ws.Range("K126").Value = 7263.599999999999
$ws.Range("L126").Value = 10949.7
$ws.Range("M126").Value = -4793.599999999999
$ws.Range("N126").Value = -15889.7

$ws.Range("H134").Value = 3759.6
$ws.Range("I134").Value = 2198.7334
$ws.Range("J134").Value = 5320.467
$ws.Range("K134").Value = 6596.2002
$ws.Range("L134").Value = 15961.401
$ws.Range("M134").Value = -4061.2002
$ws.Range("N134").Value = -21031.401

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 529.0323
$ws.Range("I5").Value = 425.66666
$ws.Range("J5").Value = 571.3182
$ws.Range("K5").Value = 1276.99998
$ws.Range("L5").Value = 1713.9546
$ws.Range("M5").Value = -1164.99998
$ws.Range("N5").Value = -1937.9546

$ws.Range("H132").Value = 842671.9399999999
$ws.Range("I132").Value = 775.8
$ws.Range("J132").Value = 1444026.2
$ws.Range("K132").Value = 6982.2
$ws.Range("L132").Value = 12996235.8
$ws.Range("M132").Value = -4452.2
$ws.Range("N132").Value = -13001295.8

$ws.Range("H135").Value = 529.0323
$ws.Range("I135").Value = 425.66666
$ws.Range("J135").Value = 571.3182
$ws.Range("K135").Value = 3830.99994
$ws.Range("L135").Value = 5141.8638
$ws.Range("M135").Value = -1295.99994
$ws.Range("N135").Value = -10211.8638

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 17384
$ws.Range("I41").Value = 2334.5715
$ws.Range("J41").Value = 70057
$ws.Range("K41").Value = 2334.5715
$ws.Range("L41").Value = 70057
$ws.Range("M41").Value = -1979.5715
$ws.Range("N41").Value = -70767

$ws.Range("H99").Value = 2200
$ws.Range("I99").Value = 2200
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 2200
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 46

$ws.Range("H102").Value = 1815.2632
$ws.Range("I102").Value = 1760.125
$ws.Range("J102").Value = 2109.3333
$ws.Range("K102").Value = 1760.125
$ws.Range("L102").Value = 2109.3333
$ws.Range("M102").Value = -138.125
$ws.Range("N102").Value = -5353.3333

$ws.Range("H113").Value = 7890.263
$ws.Range("I113").Value = 1954.4445
$ws.Range("J113").Value = 13232.5
$ws.Range("K113").Value = 1954.4445
$ws.Range("L113").Value = 13232.5
$ws.Range("M113").Value = 215.5554999999999
$ws.Range("N113").Value = -17572.5

$ws.Range("H122").Value = 2701.7646
$ws.Range("I122").Value = 2846
$ws.Range("J122").Value = 2539.5
$ws.Range("K122").Value = 8538
$ws.Range("L122").Value = 7618.5
$ws.Range("M122").Value = -6088
$ws.Range("N122").Value = -12518.5

$ws.Range("H132").Value = 4370.59
$ws.Range("I132").Value = 4558.6787
$ws.Range("J132").Value = 3891.818
$ws.Range("K132").Value = 13676.0361
$ws.Range("L132").Value = 11675.454
$ws.Range("M132").Value = -11146.0361
$ws.Range("N132").Value = -16735.454

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3545.9092
$ws.Range("I7").Value = 3235.7144
$ws.Range("J7").Value = 4088.75
$ws.Range("K7").Value = 3235.7144
$ws.Range("L7").Value = 4088.75
$ws.Range("M7").Value = -3123.7144
$ws.Range("N7").Value = -4312.75

$ws.Range("H16").Value = 1544.6666
$ws.Range("I16").Value = 1553.6
$ws.Range("J16").Value = 1500
$ws.Range("K16").Value = 1553.6
$ws.Range("L16").Value = 1500
$ws.Range("M16").Value = -1383.6
$ws.Range("N16").Value = -1840

$ws.Range("H22").Value = 22955258
$ws.Range("I22").Value = 36072120
$ws.Range("J22").Value = 753.875
$ws.Range("K22").Value = 36072120
$ws.Range("L22").Value = 753.875
$ws.Range("M22").Value = -36071825
$ws.Range("N22").Value = -1343.875

$ws.Range("H27").Value = 22955258
$ws.Range("I27").Value = 36072120
$ws.Range("J27").Value = 753.875
$ws.Range("K27").Value = 36072120
$ws.Range("L27").Value = 753.875
$ws.Range("M27").Value = -36072013
$ws.Range("N27").Value = -967.875

$ws.Range("H40").Value = 4973.591
$ws.Range("I40").Value = 4496.357
$ws.Range("J40").Value = 5808.75
$ws.Range("K40").Value = 4496.357
$ws.Range("L40").Value = 5808.75
$ws.Range("M40").Value = -4360.357
$ws.Range("N40").Value = -6080.75

$ws.Range("H46").Value = 83334530
$ws.Range("I46").Value = 125001050
$ws.Range("J46").Value = 1495.5
$ws.Range("K46").Value = 125001050
$ws.Range("L46").Value = 1495.5
$ws.Range("M46").Value = -125000862
$ws.Range("N46").Value = -1871.5

$ws.Range("H55").Value = 222.40741
$ws.Range("I55").Value = 207.42857
$ws.Range("J55").Value = 238.53847
$ws.Range("K55").Value = 207.42857
$ws.Range("L55").Value = 238.53847
$ws.Range("M55").Value = -34.42857000000001
$ws.Range("N55").Value = -584.53847

$ws.Range("H122").Value = 3202.4
$ws.Range("I122").Value = 2333
$ws.Range("J122").Value = 3575
$ws.Range("K122").Value = 6999
$ws.Range("L122").Value = 10725
$ws.Range("M122").Value = -4549
$ws.Range("N122").Value = -15625

$ws.Range("H126").Value = 3545.9092
$ws.Range("I126").Value = 3235.7144
$ws.Range("J126").Value = 4088.75
$ws.Range("K126").Value = 9707.143199999999
$ws.Range("L126").Value = 12266.25
$ws.Range("M126").Value = -7237.143199999999
$ws.Range("N126").Value = -17206.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5890117.5
$ws.Range("I62").Value = 8343058.5
$ws.Range("J62").Value = 3058.6
$ws.Range("K62").Value = 8343058.5
$ws.Range("L62").Value = 3058.6
$ws.Range("M62").Value = -8342434.5
$ws.Range("N62").Value = -4306.6

$ws.Range("H65").Value = 5890117.5
$ws.Range("I65").Value = 8343058.5
$ws.Range("J65").Value = 3058.6
$ws.Range("K65").Value = 41715292.5
$ws.Range("L65").Value = 15293
$ws.Range("M65").Value = -41712172.5
$ws.Range("N65").Value = -21533

$ws.Range("H113").Value = 495.95
$ws.Range("I113").Value = 462.16666
$ws.Range("J113").Value = 800
$ws.Range("K113").Value = 1386.49998
$ws.Range("L113").Value = 2400
$ws.Range("M113").Value = 783.5000199999999
$ws.Range("N113").Value = -6740

$ws.Range("H122").Value = 2790.0454
$ws.Range("I122").Value = 2150.9167
$ws.Range("J122").Value = 3557
$ws.Range("K122").Value = 6452.750100000001
$ws.Range("L122").Value = 10671
$ws.Range("M122").Value = -4002.750100000001
$ws.Range("N122").Value = -15571

$ws.Range("H126").Value = 1130.1177
$ws.Range("I126").Value = 1092.8
$ws.Range("J126").Value = 1183.4286
$ws.Range("K126").Value = 3278.4
$ws.Range("L126").Value = 3550.2858
$ws.Range("M126").Value = -808.3999999999996
$ws.Range("N126").Value = -8490.2858
